$wb = $excel.ActiveWorkbook

# --- Select A1 on the original (first) sheet before losing focus to the new sheet ---
[void]$wb.Worksheets.Item("Default_LoginCredentials").Range("A1").Select()

# --- Add the new "DataStructure" worksheet ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "DataStructure"

# --- Populate the new sheet with the Python-code test-case table ---
$newSheet.Range("A1").Value = "TestCaseID"
$newSheet.Range("B1").Value = "Python Code"
$newSheet.Range("C1").Value = "Expected Output"
$newSheet.Range("A2").Value = "ValidCode"
$newSheet.Range("B2").Value = 'print("hello");'
$newSheet.Range("C2").Value = "hello"

# --- Widen column C, matching the width applied on the new sheet ---
$newSheet.Columns.Item(3).ColumnWidth = 14.2

# --- Match the page setup LibreOffice applies to a freshly inserted sheet ---
$newSheet.PageSetup.LeftMargin = 56.7
$newSheet.PageSetup.RightMargin = 56.7
$newSheet.PageSetup.TopMargin = 75.8
$newSheet.PageSetup.BottomMargin = 75.8
$newSheet.PageSetup.HeaderMargin = 56.7
$newSheet.PageSetup.FooterMargin = 56.7
$newSheet.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$newSheet.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'

# --- Move the new sheet after the existing one, then make it the active tab ---
[void]$wb.Worksheets.Item("DataStructure").Move(2)
[void]$wb.Worksheets.Item("DataStructure").Range("C2").Select()
[void]$wb.Worksheets.Item("DataStructure").Activate()
